# Refresh the scraped cryptocurrency table (coinranking.com feed).
# Prices/1h-volume move with every run, and the ranking reshuffles
# slightly (row 34 'USDe' fell out of the top 50, 'Cosmos' entered
# at the bottom), so later rows shift up by one. Values are written
# directly per cell to reproduce the refreshed snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '63.806.05'
$ws.Cells.Item(2, 5).Value = '  -0.92%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.413.78'
$ws.Cells.Item(3, 5).Value = '  +0.31%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.07%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '570.72'
$ws.Cells.Item(5, 5).Value = '  -0.27%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '159.19'
$ws.Cells.Item(6, 5).Value = '  +1.03%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.05%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '3.416.16'
$ws.Cells.Item(8, 5).Value = '  +0.30%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.550'
$ws.Cells.Item(9, 5).Value = '  -9.60%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.92%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.119'
$ws.Cells.Item(11, 5).Value = '  -2.25%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.423'
$ws.Cells.Item(12, 5).Value = '  -3.67%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.006.80'
$ws.Cells.Item(13, 5).Value = '  +0.31%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +0.95%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '27.01'
$ws.Cells.Item(15, 5).Value = '  -1.80%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.0000173'
$ws.Cells.Item(16, 5).Value = '  -7.70%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '63.959.03'
$ws.Cells.Item(17, 5).Value = '  -0.70%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '3.428.14'
$ws.Cells.Item(18, 5).Value = '  +1.12%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.06'
$ws.Cells.Item(19, 5).Value = '  -4.11%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '13.56'
$ws.Cells.Item(20, 5).Value = '  -1.88%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '376.63'
$ws.Cells.Item(21, 5).Value = '  -0.59%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '7.82'
$ws.Cells.Item(22, 5).Value = '  -1.59%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.14%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '71.53'
$ws.Cells.Item(24, 5).Value = '  -0.41%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.515'
$ws.Cells.Item(25, 5).Value = '  -5.49%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.0000115'
$ws.Cells.Item(26, 5).Value = '  -1.76%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '9.58'
$ws.Cells.Item(27, 5).Value = '  -5.22%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.177'
$ws.Cells.Item(28, 5).Value = '  -0.33%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '0.999'
$ws.Cells.Item(29, 5).Value = '  -0.18%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.96'
$ws.Cells.Item(30, 5).Value = '  -2.55%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.39'
$ws.Cells.Item(31, 5).Value = '  -4.97%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.54%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '22.88'
$ws.Cells.Item(33, 5).Value = '  -1.19%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'Aptos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '7.01'
$ws.Cells.Item(34, 5).Value = '  -0.48%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.51'
$ws.Cells.Item(35, 5).Value = '  -5.19%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Monero'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '159.62'
$ws.Cells.Item(36, 5).Value = '  -0.51%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'Mantle'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.833'
$ws.Cells.Item(37, 5).Value = '  +8.43%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Stacks'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.81'
$ws.Cells.Item(38, 5).Value = '  -4.79%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Maker'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.807.28'
$ws.Cells.Item(39, 5).Value = '  -2.81%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Hedera'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0729'
$ws.Cells.Item(40, 5).Value = '  -3.23%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '25.81'
$ws.Cells.Item(41, 5).Value = '  -1.72%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'OKB'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '42.91'
$ws.Cells.Item(42, 5).Value = '  -0.04%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '6.43'
$ws.Cells.Item(43, 5).Value = '  -6.10%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Filecoin'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '4.43'
$ws.Cells.Item(44, 5).Value = '  -3.22%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '26.01'
$ws.Cells.Item(45, 5).Value = '  +1.35%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'VeChain'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0301'
$ws.Cells.Item(46, 5).Value = '  -3.82%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Bittensor'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '333.85'
$ws.Cells.Item(47, 5).Value = '  +4.92%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.37'
$ws.Cells.Item(48, 5).Value = '  +8.38%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'ONDO'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.04'
$ws.Cells.Item(49, 5).Value = '  -1.36%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Stellar'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.103'
$ws.Cells.Item(50, 5).Value = '  -4.73%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Cosmos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '6.29'
$ws.Cells.Item(51, 5).Value = '  -3.61%  '
